{"js": "// Developer Guide - Per Monitor DPI - WPF Preview : content update\n// - bump minimum OS requirement text from \"Windows 8.1 or Windows 10.\" to\n//   \"Windows 10 RS1 or higher.\" in two places (and move the \"_GoBack\"\n//   bookmark to sit right before the final \".\" of the refreshed \"Test\n//   Environment\" sentence, matching where the author's last edit landed)\n// - rename the sample variable dpiScaleInfo -> DpiScale in the\n//   RenderTargetBitmap tip\n// - drop the obsolete \"Known Issue\" paragraph about DisplayMode text\n//   rendering\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) Operating System Prerequisites paragraph.\nconst osResults = body.search(\n  \"you should be running Windows 8.1 or Windows 10.\",\n  { matchCase: true }\n);\nawait context.sync();\nif (osResults.items.length > 0) {\n  osResults.items[0].insertText(\n    \"you should be running Windows 10 RS1 or higher.\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// 2) Test Environment paragraph - update the text, then relocate the\n//    \"_GoBack\" bookmark here, right before the new trailing \".\".\nconst testEnvResults = body.search(\n  \"Run on a PC with 2 monitors, running Windows 8.1 or Windows 10.\",\n  { matchCase: true }\n);\nawait context.sync();\nif (testEnvResults.items.length > 0) {\n  testEnvResults.items[0].insertText(\n    \"Run on a PC with 2 monitors, Windows 10 RS1 or higher.\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// Move \"_GoBack\" off its old spot (end of the WinFormsHost bullet) onto the\n// point between \"...higher\" and the final \".\" we just wrote above.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst bookmarkSpotResults = body.search(\n  \"2 monitors, Windows 10 RS1 or higher\",\n  { matchCase: true }\n);\nawait context.sync();\nif (bookmarkSpotResults.items.length > 0) {\n  const justBeforePeriod = bookmarkSpotResults.items[0].getRange(\"End\");\n  justBeforePeriod.insertBookmark(\"_GoBack\");\n}\nawait context.sync();\n\n// 3) RenderTargetBitmap tip: dpiScaleInfo -> DpiScale.\nconst dpiScaleResults = body.search(\"dpiScaleInfo\", { matchCase: true });\nawait context.sync();\nif (dpiScaleResults.items.length > 0) {\n  dpiScaleResults.items[0].insertText(\"DpiScale\", \"Replace\");\n}\nawait context.sync();\n\n// 4) Remove the obsolete \"Known Issue: DisplayMode text ...\" paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Known Issue:\") === 0) {\n    paragraphs.items[i].delete();\n  }\n}\nawait context.sync();\n", "ps1": "# Developer Guide - Per Monitor DPI - WPF Preview : content update\n# - bump minimum OS requirement text from \"Windows 8.1 or Windows 10.\" to\n#   \"Windows 10 RS1 or higher.\" in two places (and move the \"_GoBack\" bookmark\n#   to sit at the end of the updated \"Test Environment\" sentence, just before\n#   the final period, matching where the author's last edit landed)\n# - rename the sample variable dpiScaleInfo -> DpiScale in the RenderTargetBitmap tip\n# - drop the obsolete \"Known Issue\" paragraph about DisplayMode text rendering\n\n$d = $word.ActiveDocument\n\n# 1) Operating System Prerequisites paragraph\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"you should be running Windows 8.1 or Windows 10.\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $rng.Text = \"you should be running Windows 10 RS1 or higher.\"\n}\n\n# 2) Test Environment paragraph - also relocate the \"_GoBack\" bookmark here,\n#    right before the trailing \".\"\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.Text = \"Run on a PC with 2 monitors, running Windows 8.1 or Windows 10.\"\n$find2.MatchCase = $true\nif ($find2.Execute()) {\n    $sentenceStart = $rng2.Start\n    $newSentence = \"Run on a PC with 2 monitors, Windows 10 RS1 or higher.\"\n    $rng2.Text = $newSentence\n    $sentenceEnd = $sentenceStart + $newSentence.Length\n\n    # Move the \"_GoBack\" bookmark from its old spot (end of the WinFormsHost\n    # bullet) to just before the final \".\" of this sentence.\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n    $bmPoint = $d.Range($sentenceEnd - 1, $sentenceEnd - 1)\n    $d.Bookmarks.Add(\"_GoBack\", $bmPoint) | Out-Null\n}\n\n# 3) RenderTargetBitmap tip: dpiScaleInfo -> DpiScale\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.Text = \"dpiScaleInfo\"\n$find3.MatchCase = $true\nif ($find3.Execute()) {\n    $rng3.Text = \"DpiScale\"\n}\n\n# 4) Remove the obsolete \"Known Issue: DisplayMode text ...\" paragraph\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"Known Issue:*DisplayMode text*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
